$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
$ws.Range("D7").Value = 43373
